$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.528.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "'1.865.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "'0.9971"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'246.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'0.7050"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("D7").Value = "'0.9976"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.07783"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.3092"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'23.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("D11").Value = "'0.07842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "'5.186"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "'93.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("D14").Value = "'1.854.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "'0.6997"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.74%  "
$ws.Range("D16").Value = "'6.660"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "'0.000008397"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "'29.457.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "'244.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "'2.100.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'12.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'0.9987"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'7.603"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").Value = "'0.9988"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'0.1526"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "
$ws.Range("D26").Value = "'8.959"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").Value = "'159.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "'18.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'1.543"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "'4.267"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").Value = "'4.231"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").Value = "'1.203"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "'0.05163"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D34").Value = "'0.7983"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").Value = "'1.936"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").Value = "'1.157"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.11%  "
$ws.Range("D37").Value = "'2.693"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "
$ws.Range("D38").Value = "'1.335.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.95%  "
$ws.Range("D39").Value = "'0.01881"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "'2.732"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "'0.9609"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.69%  "
$ws.Range("D42").Value = "'6.096"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.50%  "
$ws.Range("D43").Value = "'107.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.27%  "
$ws.Range("D44").Value = "'0.9984"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'9.829"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.85%  "
$ws.Range("D46").Value = "'2.002.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'65.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").Value = "'0.5199"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'1.787"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").Value = "'7.044"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.00%  "
